$d = $word.ActiveDocument

# --- Swap "Post-Conditions" <-> "Preconditions" across the whole document ---
# A direct two-way replace would collide (the new "Preconditions" text created by
# replacing "Post-Conditions" would itself get matched by the second replace).
# Route through a unique placeholder to perform a clean swap.
$d.Content.Find.Execute("Post-Conditions", $true, $true, $false, $false, $false, $true, 1, $false, "@@TEMP_SWAP_PLACEHOLDER@@", 2)
$d.Content.Find.Execute("Preconditions", $true, $true, $false, $false, $false, $true, 1, $false, "Post-conditions", 2)
$d.Content.Find.Execute("@@TEMP_SWAP_PLACEHOLDER@@", $true, $true, $false, $false, $false, $true, 1, $false, "Preconditions", 2)

# --- Normalize a few runs that get merged in the reference edit ---
# (Touching text that spans these run boundaries merges the runs into one,
# matching the authoring tool's normal behavior when text near a split is edited.)
$mergeTargets = @(
    "Sys-admin-05",
    "Sys-admin-06",
    "Sys-admin-07",
    "Sys-admin-08",
    "Generate reports about teachers",
    "Generate reports about organizations",
    "Generate reports about admins",
    "Generate reports about courses"
)
foreach ($t in $mergeTargets) {
    $d.Content.Find.Execute($t, $true, $true, $false, $false, $false, $true, 1, $false, $t, 2)
}

# --- Move the "_GoBack" bookmark ---
# It previously lived in the empty trailing paragraph after the last table; it now
# belongs right after the text in the last table's "Post-conditions" label cell.
$tbl = $d.Tables.Item($d.Tables.Count)
$cell = $tbl.Cell(7, 1)
$para = $cell.Range.Paragraphs.Item(1)
$endPos = $para.Range.End - 1

# Creating a bookmark at a truly collapsed position right at the end of a run can
# misplace it, so temporarily insert a marker character, anchor the bookmark next
# to it, then remove the marker, leaving the bookmark collapsed in the right spot.
$insertRange = $d.Range($endPos, $endPos)
$insertRange.InsertAfter("X")
$target = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $target)
$markerRange = $d.Range($endPos, $endPos + 1)
$markerRange.Delete()
